$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: the Price column (D) holds values that look numeric (e.g. "1.004",
# "322.54") but must stay as literal text, matching the source data which
# is stored as inline strings. Assigning such a string straight to
# .Value lets Excel's COM layer auto-convert it to a floating point
# number, so we force the target cell's NumberFormat to Text ("@")
# immediately before writing values that would otherwise parse as a
# plain number. Multi-dot values (e.g. "27.783.43") and the Volume column
# (which always carries padding spaces/percent signs) are never
# re-interpreted as numbers by Excel, so no extra step is required there.

# --- Reorder rows 37-39: the coin that was in row 39 (Hedera) moves to
#     row 37, the coin in row 37 (Algorand) moves to row 38, and the coin
#     in row 38 (InternetComputer(DFINITY)) moves to row 39. Column A
#     (the index number) stays put; only B (Coin), C (Link) and D (Price)
#     change identity/value, and E (Volume) gets the new figures below.
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06076"
$ws.Range("E37").Value = "  -0.11%  "

$ws.Range("B38").Value = "Algorand"
$ws.Range("C38").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2108"
$ws.Range("E38").Value = "  -1.96%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.041"
$ws.Range("E39").Value = "  +0.15%  "

# --- Update Price (D) and Volume(1h) (E) columns for each remaining row
# Row 2
$ws.Range("D2").Value = "27.783.43"
$ws.Range("E2").Value = "  +1.19%  "
# Row 3
$ws.Range("D3").Value = "1.761.67"
$ws.Range("E3").Value = "  +0.46%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  -0.38%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "322.54"
$ws.Range("E5").Value = "  -0.61%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  -0.15%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4248"
$ws.Range("E7").Value = "  -3.86%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3600"
$ws.Range("E8").Value = "  -2.39%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.35"
$ws.Range("E9").Value = "  -1.06%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07436"
$ws.Range("E10").Value = "  -2.71%  "
# Row 11
$ws.Range("E11").Value = "  -0.80%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9996"
$ws.Range("E12").Value = "  -0.40%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.36"
$ws.Range("E13").Value = "  -0.90%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.103"
$ws.Range("E14").Value = "  -0.57%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.286"
$ws.Range("E15").Value = "  -1.81%  "
# Row 16
$ws.Range("D16").Value = "1.799.45"
$ws.Range("E16").Value = "  +2.00%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "90.90"
$ws.Range("E17").Value = "  +1.12%  "
# Row 18
$ws.Range("E18").Value = "  -0.73%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06390"
$ws.Range("E19").Value = "  +2.19%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9991"
$ws.Range("E20").Value = "  -0.37%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.14"
$ws.Range("E21").Value = "  -1.19%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.950"
$ws.Range("E22").Value = "  -3.29%  "
# Row 23
$ws.Range("D23").Value = "27.827.50"
$ws.Range("E23").Value = "  +1.12%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.30"
$ws.Range("E24").Value = "  -1.82%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.151"
$ws.Range("E25").Value = "  -7.07%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.88"
$ws.Range("E26").Value = "  +5.42%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.22"
$ws.Range("E27").Value = "  -1.27%  "
# Row 28
$ws.Range("D28").Value = "1.996.56"
$ws.Range("E28").Value = "  +1.97%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.145"
$ws.Range("E29").Value = "  -6.09%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.65"
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.167"
$ws.Range("E31").Value = "  -0.98%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.658"
$ws.Range("E32").Value = "  -0.77%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09024"
$ws.Range("E33").Value = "  -1.72%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.515"
$ws.Range("E34").Value = "  -3.20%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.64"
$ws.Range("E35").Value = "  +0.58%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02311"
$ws.Range("E36").Value = "  +0.20%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6403"
$ws.Range("E40").Value = "  -0.40%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.184"
$ws.Range("E41").Value = "  +1.31%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  -0.18%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.890"
$ws.Range("E43").Value = "  -0.70%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.394"
$ws.Range("E44").Value = "  +0.01%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.67"
$ws.Range("E45").Value = "  -0.34%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5955"
$ws.Range("E46").Value = "  +0.20%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.705"
$ws.Range("E47").Value = "  -0.40%  "
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "124.39"
$ws.Range("E48").Value = "  -1.26%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.982"
$ws.Range("E49").Value = "  -0.17%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.147"
$ws.Range("E50").Value = "  +1.60%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06878"
$ws.Range("E51").Value = "  -0.07%  "
